$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new rows to the bottom of the Sheet1 table.
$ws.Range("A11").Value = "催账单"
$ws.Range("B11").Value = "【标题】Request for 月份 Invoice`n【正文】Hi XXX Team,`nHope this email finds you well.`nCould you please kindly provide the shipping fee invoice for 月份? We need it for accounting and settlement purposes.`nIf there is any additional information required from our side, please feel free to let us know.`nThank you for your assistance.`nBest regards,"

$ws.Range("A12").Value = "账单申诉-重量差异"
$ws.Range("B12").Value = "We hope this email finds you well.`nPlease disregard our previous message. Kindly refer to the attached document for the correct information: we need your assistance to verify the weight discrepancies involving XX orders, with a total amount of XX.`nThank you for your prompt attention to these matters. We look forward to your reply.`nBest regards,"

# Excel auto-expands row height when multi-line text is entered into a
# wrapped cell; put the new rows back to the sheet's default (un-customized)
# height so the saved XML matches a plain data row like the others.
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()

# Update selection / view state to match the authored change.
$ws.Range("K24").Select() | Out-Null
